$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that were missing the "White (%)" value in column O because the
# "Native Hawaiian (%)" value in column L had been omitted, shifting every
# later race-percentage column one place to the left. Re-insert the missing
# L value (0) and shift the existing L/M/N values right into M/N/O.
$rows = 21, 24, 47, 50

foreach ($r in $rows) {
    $colO = $ws.Cells.Item($r, 15)  # O
    $colN = $ws.Cells.Item($r, 14)  # N
    $colM = $ws.Cells.Item($r, 13)  # M
    $colL = $ws.Cells.Item($r, 12)  # L

    # Shift existing values one column to the right, starting from the
    # rightmost column so we don't clobber values before they're copied.
    $colO.Value2 = $colN.Value2
    $colN.Value2 = $colM.Value2
    $colM.Value2 = $colL.Value2
    $colL.Value2 = 0
}

# Update the saved selection to match the final interactive state.
$ws.Range("M52").Select()
